$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns at I:J ("Owner" / "Owner Email"). Everything that
#    used to live at column I onward (Gross Area, GHG Intensity, ... ) shifts
#    right by two columns, which this engine also re-points cell formulas
#    for automatically.
# ---------------------------------------------------------------------------
$ws.Range("I1:J1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. Header row updates
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Audit Template Building ID"
$ws.Range("B1").Value = "Portfolio Manager Building ID"

# ---------------------------------------------------------------------------
# 3. Column B (old numeric "Building ID") becomes the textual
#    "Portfolio Manager Building ID" values. Force text storage so the
#    value is kept as a shared string (matching the source data) instead of
#    being reinterpreted as a number.
#    Row 2 already carries the sheet's quotePrefix style (s="2"), so a
#    leading apostrophe keeps it on that same style with no new style rows.
#    Rows 3-10 have no explicit style, so a temporary Text number format is
#    used to force string storage and then cleared back to Normal.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 2).Value = "'21537666"

$pmIds = @("21537667","21537668","21537669","21537670","21537671","21537672","21537673","21537674")
for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 3
    $cell = $ws.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $pmIds[$i]
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 4. New Owner / Owner Email columns
# ---------------------------------------------------------------------------
$ws.Range("I1").Value = "Owner"
$ws.Range("J1").Value = "Owner Email"

$owners = @("Company A","Company B","Company C","Company D","Company E","Company F","Company G","Company H","Company I")
$emails = @("admin@companya.com","admin@companyb.com","admin@companyc.com","admin@companyd.com","admin@companye.com","admin@companyf.com","admin@companyg.com","admin@companyh.com","admin@companyi.com")

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $owners[$i]
    $ws.Cells.Item($row, 10).Value = $emails[$i]
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 10), "mailto:" + $emails[$i]) | Out-Null
}

# ---------------------------------------------------------------------------
# 5. Conditional formatting referenced the pre-insert column letters; this
#    engine does not shift those ranges/formulas automatically the way real
#    Excel does, so rebuild the rules against the post-insert layout.
# ---------------------------------------------------------------------------
$ws.Cells.FormatConditions.Delete()

# was I2:I10 -> now K2:K10  (Gross Area)
$fc = $ws.Range("K2:K10").FormatConditions.Add(2, 0, '=$P2="Yes"')
$fc.StopIfTrue = $true
$fc = $ws.Range("K2:K10").FormatConditions.Add(1, 2, "1000", "1000000")
$fc = $ws.Range("K2:K10").FormatConditions.Add(2, 0, '=$R2="Yes"')
$fc.StopIfTrue = $true

# was K2:K10 -> now M2:M10 (Site EUI)
$fc = $ws.Range("M2:M10").FormatConditions.Add(1, 2, "40", "1000")

# was O9:O10 -> now Q9:Q10 (GHGI Target Year)
$fc = $ws.Range("Q9:Q10").FormatConditions.Add(2, 0, '=$Q8="Yes"')
$fc.StopIfTrue = $true
$fc = $ws.Range("Q9:Q10").FormatConditions.Add(2, 0, '=$T8="Yes"')
$fc.StopIfTrue = $true

# was N7:O8 N9:N10 -> now P7:Q8 P9:P10
$fc = $ws.Range("P7:Q8,P9:P10").FormatConditions.Add(2, 0, '=$S6="Yes"')
$fc.StopIfTrue = $true
$fc = $ws.Range("P7:Q8,P9:P10").FormatConditions.Add(2, 0, '=$T6="Yes"')
$fc.StopIfTrue = $true

# was P2:P10 -> now R2:R10 (EUI Target)
$fc = $ws.Range("R2:R10").FormatConditions.Add(1, 2, "40", "1000")

# was N2:O6 O2:Q10 -> now P2:Q6 Q2:S10
$fc = $ws.Range("P2:Q6,Q2:S10").FormatConditions.Add(2, 0, '=$S2="Yes"')
$fc.StopIfTrue = $true
$fc = $ws.Range("P2:Q6,Q2:S10").FormatConditions.Add(2, 0, '=$T2="Yes"')
$fc.StopIfTrue = $true

# was J2:K10 -> now L2:M10 (Total GHG Emissions Intensity / Site EUI)
$fc = $ws.Range("L2:M10").FormatConditions.Add(2, 0, '=$P2="Yes"')
$fc.StopIfTrue = $true
$fc = $ws.Range("L2:M10").FormatConditions.Add(2, 0, '=$R2="Yes"')
$fc.StopIfTrue = $true

# ---------------------------------------------------------------------------
# 6. Selection cosmetics to mirror the column-insert side effect of the
#    authoring session (selecting the freshly inserted columns).
# ---------------------------------------------------------------------------
$ws.Range("I1:J1048576").Select()
